# Update dashboards - 2026-02-27
# Weekly refresh of the "Latest Period" dates and the Present/Lag1-4 history
# columns for the daily/weekly interest-rate rows (T5YIFR, T10YIE, DFF, DGS2,
# DGS5, DGS10, MORTGAGE30US, DBAA), plus the matching "new data this week"
# highlight move (the yellow highlight follows the row whose date advanced
# the most - row 51, MORTGAGE30US).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 (5Y Forward Inflation Expectation Rate, T5YIFR) ---
$ws.Range("N29").Value = 46079
$ws.Range("Q29").Value = 2.13
$ws.Range("R29").Value = 2.14
$ws.Range("T29").Value = 2.12
$ws.Range("U29").Value = 2.13

# --- Row 30 (10Y Breakeven Inflation Rate, T10YIE) ---
$ws.Range("N30").Value = 46079
$ws.Range("R30").Value = 2.28
$ws.Range("T30").Value = 2.26
$ws.Range("U30").Value = 2.28

# --- Row 47 (Fed Funds Rate, DFF) ---
$ws.Range("N47").Value = 46078

# --- Row 48 (2-Year Treasury, DGS2) ---
$ws.Range("N48").Value = 46078
$ws.Range("Q48").Value = 3.45
$ws.Range("S48").Value = 3.43
$ws.Range("T48").Value = 3.48

# --- Row 49 (5-Year Treasury, DGS5) ---
$ws.Range("N49").Value = 46078
$ws.Range("R49").Value = 3.61
$ws.Range("S49").Value = 3.59
$ws.Range("U49").Value = 3.65

# --- Row 50 (10-Year Treasury, DGS10) ---
$ws.Range("N50").Value = 46078
$ws.Range("Q50").Value = 4.05
$ws.Range("R50").Value = 4.04
$ws.Range("S50").Value = 4.03
$ws.Range("U50").Value = 4.08

# --- Row 51 (30-Year Mortgage Rate, MORTGAGE30US) ---
$ws.Range("N51").Value = 46076
$ws.Range("Q51").Value = 5.98
$ws.Range("R51").Value = 6.01
$ws.Range("S51").Value = 6.09
$ws.Range("T51").Value = 6.11
$ws.Range("U51").Value = 6.1

# Row 51's "Latest Period" date jumped forward the most (new data just in),
# so the "updated this week" yellow highlight moves onto it - copy the
# existing highlighted date-cell format (e.g. N13) over rather than inventing
# a new fill.
$ws.Range("N13").Copy()
$ws.Range("N51").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 52 (BAA Corporate Bond Yield, DBAA) ---
$ws.Range("N52").Value = 46078
$ws.Range("Q52").Value = 5.78
$ws.Range("R52").Value = 5.77
$ws.Range("S52").Value = 5.76
$ws.Range("T52").Value = 5.77
